$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("relays")

# Update row 3 values (C3, D3, E3)
$ws.Range("C3").Value = 5
$ws.Range("D3").Value = 12
$ws.Range("E3").Value = 0.5

# Update row 4 value (H4)
$ws.Range("H4").Value = 0

# Update the active selection to H4
$ws.Range("H4").Select()
